# This workbook lists weekly price records (rows 2-16) for Níspero at the
# "Vega Modelo de Temuco" market. The edit re-shuffles the per-row data
# (Fecha, Calidad, Volumen, Precio minimo/maximo/promedio, Unidad de
# comercializacion, Origen, Precio $/Kg, Kg/unidad) across rows 2-13,15,16
# while row 14 and the remaining descriptive columns (A,B,C,E,F,G,H,I,J,K)
# stay untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per row for columns: D(4) L(12) M(13) N(14) O(15) P(16) Q(17) R(18) S(19) T(20)
$targets = @(
    @{ Row=2;  D=44166; L="Segunda"; M=20;  N=12000; O=12000; P=12000; Q="$/caja 18 kilos";   R="La Ligua";               S=667;  T=18 },
    @{ Row=3;  D=44858; L="Primera"; M=90;  N=20000; O=20000; P=20000; Q="$/bandeja 5 kilos";  R="Provincia de Quillota";  S=4000; T=5  },
    @{ Row=4;  D=44515; L="Primera"; M=80;  N=28000; O=28000; P=28000; Q="$/bandeja 10 kilos"; R="Provincia de Los Andes"; S=2800; T=10 },
    @{ Row=5;  D=44519; L="Primera"; M=30;  N=28000; O=28000; P=28000; Q="$/bandeja 10 kilos"; R="Provincia de Quillota";  S=2800; T=10 },
    @{ Row=6;  D=44511; L="Primera"; M=45;  N=28000; O=28000; P=28000; Q="$/bandeja 10 kilos"; R="Provincia de Los Andes"; S=2800; T=10 },
    @{ Row=7;  D=44511; L="Primera"; M=45;  N=3200;  O=3200;  P=3200;  Q="$/bandeja 10 kilos"; R="Provincia de Quillota";  S=320;  T=10 },
    @{ Row=8;  D=44859; L="Primera"; M=30;  N=20000; O=20000; P=20000; Q="$/bandeja 5 kilos";  R="Provincia de Quillota";  S=4000; T=5  },
    @{ Row=9;  D=44879; L="Primera"; M=25;  N=30000; O=30000; P=30000; Q="$/bandeja 10 kilos"; R="Provincia de Quillota";  S=3000; T=10 },
    @{ Row=10; D=44868; L="Primera"; M=30;  N=14000; O=14000; P=14000; Q="$/bandeja 5 kilos";  R="Provincia de Quillota";  S=2800; T=5  },
    @{ Row=11; D=44483; L="Primera"; M=35;  N=10000; O=10000; P=10000; Q="$/bandeja 5 kilos";  R="Provincia de Quillota";  S=2000; T=5  },
    @{ Row=12; D=44874; L="Primera"; M=40;  N=25000; O=25000; P=25000; Q="$/bandeja 10 kilos"; R="Provincia de Quillota";  S=2500; T=10 },
    @{ Row=13; D=44466; L="Primera"; M=80;  N=11000; O=11000; P=11000; Q="$/bandeja 5 kilos";  R="La Ligua";               S=2200; T=5  },
    @{ Row=15; D=44503; L="Primera"; M=50;  N=28000; O=28000; P=28000; Q="$/bandeja 10 kilos"; R="Provincia de Quillota";  S=2800; T=10 },
    @{ Row=16; D=44496; L="Primera"; M=55;  N=28000; O=28000; P=28000; Q="$/bandeja 10 kilos"; R="Provincia de Quillota";  S=2800; T=10 }
)

foreach ($t in $targets) {
    $r = $t.Row
    $ws.Cells.Item($r, 4).Value = $t.D
    $ws.Cells.Item($r, 12).Value = $t.L
    $ws.Cells.Item($r, 13).Value = $t.M
    $ws.Cells.Item($r, 14).Value = $t.N
    $ws.Cells.Item($r, 15).Value = $t.O
    $ws.Cells.Item($r, 16).Value = $t.P
    $ws.Cells.Item($r, 17).Value = $t.Q
    $ws.Cells.Item($r, 18).Value = $t.R
    $ws.Cells.Item($r, 19).Value = $t.S
    $ws.Cells.Item($r, 20).Value = $t.T
}
